$d = $word.ActiveDocument

# --- 1. Merge the split runs of the opening sentence into a single run. ---
# (Text is unchanged - only the run structure is normalized, matching the
#  author's edit which collapsed four adjacent runs into one.)
$d.Content.Find.Execute(
    "I arrive at the music school an hour before I was supposed to, finding that the auditorium hasn’t even been opened for the audience yet. Prim’s been rubbing off on me, huh.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I arrive at the music school an hour before I was supposed to, finding that the auditorium hasn’t even been opened for the audience yet. Prim’s been rubbing off on me, huh.",
    2) | Out-Null

# --- 2. Remove the "?Iris: Are you..." / "?Iris: Pro?" / "Pro: Yeah, that's me."
#        paragraphs entirely. ---
$d.Content.Find.Execute(
    "?Iris: Are you…^p?Iris: Pro?^pPro: Yeah, that’s me.^p",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2) | Out-Null

# --- 3. Update the "I blink twice..." line describing the college student. ---
$d.Content.Find.Execute(
    "She looks like she’s a college student, and her hand is wrapped in a cast…",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A college student, with a hand wrapped in a cast…",
    2) | Out-Null

# --- 4. "Pro: You're Prim's sister?" -> "Pro: Prim's sister, uh..." ---
$d.Content.Find.Execute(
    "Pro: You’re Prim’s sister?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pro: Prim’s sister, uh…",
    2) | Out-Null

# --- 5. Remove "She nods, and I let out an internal sigh of relief." paragraph. ---
$d.Content.Find.Execute(
    "She nods, and I let out an internal sigh of relief.^p",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2) | Out-Null

# --- 6. "Iris (neutral smiling): I'm Iris." -> "Iris (neutral disappointed): Iris." ---
$d.Content.Find.Execute(
    "Iris (neutral smiling): I’m Iris.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Iris (neutral disappointed): Iris.",
    2) | Out-Null

# --- 7. "Iris (neutral curious): Has Prim told you about me?" -> "Pro: Right. Sorry." ---
$d.Content.Find.Execute(
    "Iris (neutral curious): Has Prim told you about me?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Pro: Right. Sorry.",
    2) | Out-Null

# --- 8. Remove the standalone "Pro: Yeah." paragraph. ---
$d.Content.Find.Execute(
    "Pro: Yeah.^p",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2) | Out-Null

# --- 9. "Iris (neutral neutral): I see." -> "Iris (neutral neutral): Don't worry about it." ---
$d.Content.Find.Execute(
    "Iris (neutral neutral): I see.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Iris (neutral neutral): Don’t worry about it.",
    2) | Out-Null

# --- 10. Insert a new paragraph right before "To my surprise, she stands up
#         straight and bows." describing Iris scrutinizing Pro's outfit. ---
$r = $d.Content
$r.Find.Execute(
    "To my surprise, she stands up straight and bows.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    0) | Out-Null
$r.InsertBefore("She stares at me as if scrutinizing my entire being, her eyes glancing over the casual t-shirt and jeans I threw on before I left. Now that I think about it, aren’t concerts supposed to be formal events? And isn’t Iris a professional pianist, so maybe she disapproves…?`r")

# --- 11. Reword that paragraph itself. ---
$d.Content.Find.Execute(
    "To my surprise, she stands up straight and bows.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "But to my surprise, instead of saying anything she stands up straight and bows.",
    2) | Out-Null

# --- 12. Merge the split runs of the "eavesdropping" sentence into one run. ---
$d.Content.Find.Execute(
    "I instinctively freeze up, deciding not to tell her that I was there eavesdropping.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I instinctively freeze up, deciding not to tell her that I was there eavesdropping.",
    2) | Out-Null
